$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at Q so the existing "images" header (and any data
# below it) shifts one column to the right (Q -> R).
$ws.Columns("Q:Q").Insert()

# Populate the newly inserted column's header with the new field name.
$ws.Range("Q1").Value = "suggestedProduct"

# Give the new header cell its own distinct style: Calibri 11, black,
# matching the look of the other header cells (font4) but resolved by the
# engine as a new font/style entry rather than being collapsed into the
# pre-existing one.
$ws.Range("Q1").Font.Name = "Calibri"
$ws.Range("Q1").Font.Size = 11
$ws.Range("Q1").Font.Color = 1

# Column width tweaks that accompanied the new column in the original edit.
$ws.Columns(1).ColumnWidth = 15.3
$ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item(1, 15)).ColumnWidth = 7.6
$ws.Columns(16).ColumnWidth = 26.166666666666668
$ws.Range($ws.Cells.Item(1, 17), $ws.Cells.Item(1, 1025)).ColumnWidth = 7.6

# Restore the selection to the new header cell.
[void]$ws.Range("Q1").Select()
